# Generate Report for Handoff
#
# Author's commit adds a new localization handoff record for the file
# "4766ce3c-ee45-4edc-835c-e0046d90834d...md" as a new row (row 3) in the
# Overview sheet and in each language sheet (zh-cn, de-de), and grows the
# backing Excel Tables + AutoFilter ranges + sheet dimensions to match.

$wb = $excel.ActiveWorkbook

$newMdName    = "4766ce3c-ee45-4edc-835c-e0046d90834dooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newMdDisplay = "e2e\4766ce3c-ee45-4edc-835c-e0046d90834dooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newXlfZh     = "4766ce3c-ee45-4edc-835c-e0046d90834doooooooooooooooooooooooooooooooooooooooo.2105cd5bac64ab052b0ce9ed8cd57175846b00eb.zh-cn.xlf"
$newXlfDe     = "4766ce3c-ee45-4edc-835c-e0046d90834doooooooooooooooooooooooooooooooooooooooo.2105cd5bac64ab052b0ce9ed8cd57175846b00eb.de-de.xlf"
$readyStatus  = "Ready for handoff"
$commitHash   = "2105cd5bac64ab052b0ce9ed8cd57175846b00eb"
$repoBlobBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/"

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1 / table "Overview")
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item("Overview")
$tblOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newMdName
$wsOverview.Range("B3").Value = $newMdDisplay
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), ($repoBlobBase + $newMdName), "", "", $newMdDisplay) | Out-Null
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = $readyStatus
$wsOverview.Range("F3").Value = $readyStatus
$wsOverview.Range("G3").Value = "2016-09-01 14:35:44"

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table "zh-cn")
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$tblZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$tblZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $newMdName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), ($repoBlobBase + $newMdName), "", "", $newMdDisplay) | Out-Null
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $readyStatus
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $newXlfZh
$wsZhCn.Range("H3").Value = "2016-09-01 14:35:39"
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3 / table "de-de")
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$tblDeDe = $wsDeDe.ListObjects.Item("de-de")
$tblDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $newMdName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), ($repoBlobBase + $newMdName), "", "", $newMdDisplay) | Out-Null
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $readyStatus
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $newXlfDe
$wsDeDe.Range("H3").Value = "2016-09-01 14:35:44"
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

# ---------------------------------------------------------------------
# Column width tweaks carried by the same commit (widen the "Latest HO
# Xliff Generate Date" / HO-datetime columns to fit the new timestamps).
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
